{"js": "// Adi\u00e7\u00e3o de aviso aos estudantes\n// Insert a bold/italic/red warning paragraph right after the standalone\n// \"github.com\" paragraph near the end of the document.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the paragraph whose *entire* text is exactly \"github.com\" (not the\n// \"https://desktop.github.com/\" hyperlink paragraph earlier in the doc,\n// which merely contains that text as a substring).\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"github.com\") {\n    targetParagraph = paragraphs.items[i];\n  }\n}\n\nif (!targetParagraph) {\n  throw new Error('Could not find the \"github.com\" paragraph to anchor the new warning after.');\n}\n\nconst warningText = \"VOC\u00ca USAR\u00c1 ESSAS FERRAMENTAS NO CASE E AULA E NO CHALLENGE!\";\nconst newParagraph = targetParagraph.insertParagraph(warningText, Word.InsertLocation.after);\n\n// The source paragraph carried a first-line indent that the new warning\n// paragraph should not have.\nnewParagraph.paragraphFormat.firstLineIndent = 0;\nawait context.sync();\n\n// Match the formatting from the diff: bold, italic, red, no underline\n// (the new paragraph already inherited the \"Hyperlink\" character style from\n// the paragraph it was split off from). Apply to the run text itself...\nconst contentRange = newParagraph.getRange(Word.RangeLocation.content);\ncontentRange.font.bold = true;\ncontentRange.font.boldBidirectional = true;\ncontentRange.font.italic = true;\ncontentRange.font.italicBidirectional = true;\ncontentRange.font.color = \"#FF0000\";\ncontentRange.font.underline = Word.UnderlineType.none;\nawait context.sync();\n\n// ...and to the paragraph mark itself, so the paragraph-level rPr (the\n// formatting that would apply to whatever is typed next at the end of this\n// paragraph) also reflects the same bold/italic/red/no-underline look.\nconst endRange = newParagraph.getRange(Word.RangeLocation.end);\nendRange.font.bold = true;\nendRange.font.italic = true;\nendRange.font.color = \"#FF0000\";\nendRange.font.underline = Word.UnderlineType.none;\nawait context.sync();\n", "ps1": "# Adi\u00e7\u00e3o de aviso aos estudantes\n# Insert a bold/italic/red warning paragraph right after the standalone\n# \"github.com\" paragraph near the end of the document.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph whose *entire* text is exactly \"github.com\" (not the\n# \"https://desktop.github.com/\" hyperlink paragraph earlier in the doc,\n# which merely contains that text as a substring).\n$targetIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    $trimmed = $t.TrimEnd([char]13, [char]7)\n    if ($trimmed -eq \"github.com\") {\n        $targetIndex = $i\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw 'Could not find the \"github.com\" paragraph to anchor the new warning after.'\n}\n\n$targetPara = $d.Paragraphs.Item($targetIndex)\n$endOfTarget = $targetPara.Range\n$endOfTarget.Collapse(0)  # wdCollapseEnd\n$endOfTarget.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Item($targetIndex + 1)\n$newPara.Range.Text = \"VOC\u00ca USAR\u00c1 ESSAS FERRAMENTAS NO CASE E AULA E NO CHALLENGE!\"\n\n# Match the formatting from the diff: bold, italic, red, no underline\n# (the new paragraph already inherited the \"Hyperlink\" character style from\n# the paragraph it was split off from).\n$newRng = $newPara.Range\n$newRng.Font.Bold = 1\n$newRng.Font.BoldBi = 1\n$newRng.Font.Italic = 1\n$newRng.Font.ItalicBi = 1\n$newRng.Font.Color = 255\n$newRng.Font.Underline = 0\n\n# The source paragraph carried a first-line indent that the new warning\n# paragraph should not have.\n$newPara.Format.FirstLineIndent = 0\n"}
